$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header for column BB (day after BA1 = 43973)
$ws.Range("BB1").Value = 43974

# Daily totals for column BB, mirroring column BA's per-row structure
$ws.Range("BB2").Value = 109
$ws.Range("BB3").Value = 80
$ws.Range("BB4").Value = 42
$ws.Range("BB5").Value = 31
$ws.Range("BB6").Value = 5
$ws.Range("BB7").Value = 2407
$ws.Range("BB8").Value = 23
$ws.Range("BB9").Value = 715
$ws.Range("BB10").Value = 0
$ws.Range("BB11").Value = 12
$ws.Range("BB12").Value = 1
$ws.Range("BB13").Value = 9
$ws.Range("BB14").Value = 12
$ws.Range("BB15").Value = 2
$ws.Range("BB16").Value = 11
$ws.Range("BB17").Value = 35
$ws.Range("BB18").Value = 45
$ws.Range("BB19").Value = 204

# Totals row formula, same pattern as BA20
$ws.Range("BB20").Formula = "=SUM(BB2:BB19)"

# Update the active selection to match the newly populated column
$ws.Range("BB2:BB20").Select()
